# 自动更新Excel文件 - 2026-01-13 23:14:05
# 每日巡检：“剩余”天数（E列）递减1；当“剩余”归零（原值为1）时，
# 说明本轮周期结束，重置为10，并将“开始时间”（F列）顺延10天，开启下一轮。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    # 跳过空行或数据异常（开始时间不是合法的8位日期）的行
    if ($null -eq $eVal -or $eVal -eq "") {
        continue
    }
    if ($null -eq $fVal -or $fVal -eq "") {
        continue
    }

    $fStr = [string]([int64]$fVal)
    if ($fStr.Length -ne 8) {
        continue
    }

    $eNum = [int]$eVal

    if ($eNum -eq 1) {
        $eCell.Value = 10
        $fCell.Value = [int64]$fVal + 10
    } else {
        $eCell.Value = $eNum - 1
    }
}
